# Fruta / hortaliza, semanal
# Insert two new weekly-price rows at the top of the "Piña" data block (row 307),
# pushing all existing rows from 307 onward down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 307 (each Insert() pushes the current
# row 307 and everything below it down by one row).
$ws.Rows.Item(307).Insert()
$ws.Rows.Item(307).Insert()

# ---- New row 307 ----
$ws.Range("A307").Value = 7
$ws.Range("B307").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C307").Value = "Ñuble"
$ws.Range("D307").Value = 45127
$ws.Range("E307").Value = 16
$ws.Range("F307").Value = "Fruta"
$ws.Range("G307").Value = 100108
$ws.Range("H307").Value = "Tropicales y subtropicales"
$ws.Range("I307").Value = 100108005
$ws.Range("J307").Value = "Piña"
$ws.Range("K307").Value = "Caramelo"
$ws.Range("L307").Value = "Primera"
$ws.Range("M307").Value = 80
$ws.Range("N307").Value = 24000
$ws.Range("O307").Value = 24000
$ws.Range("P307").Value = 24000
$ws.Range("Q307").Value = "$/caja 10 unidades"
$ws.Range("R307").Value = "Ecuador"
$ws.Range("S307").Value = 2400
$ws.Range("T307").Value = 10

# ---- New row 308 ----
$ws.Range("A308").Value = 7
$ws.Range("B308").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C308").Value = "Ñuble"
$ws.Range("D308").Value = 45127
$ws.Range("E308").Value = 16
$ws.Range("F308").Value = "Fruta"
$ws.Range("G308").Value = 100108
$ws.Range("H308").Value = "Tropicales y subtropicales"
$ws.Range("I308").Value = 100108005
$ws.Range("J308").Value = "Piña"
$ws.Range("K308").Value = "Caramelo"
$ws.Range("L308").Value = "Primera"
$ws.Range("M308").Value = 30
$ws.Range("N308").Value = 23000
$ws.Range("O308").Value = 23000
$ws.Range("P308").Value = 23000
$ws.Range("Q308").Value = "$/caja 12 unidades"
$ws.Range("R308").Value = "Ecuador"
$ws.Range("S308").Value = 1917
$ws.Range("T308").Value = 12

# Make sure the D column on the two new rows keeps the date number format
# used throughout the column (style index 2 in the original workbook).
$ws.Range("D307:D308").NumberFormat = "YYYY-MM-DD HH:MM:SS"
